$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column values are stored as plain text in this sheet (e.g. multi-dot
# figures like "27.457.58"). For cells whose new text would otherwise look like a
# plain number to Excel's smart-entry parser, force the cell to Text format first
# so the literal string is preserved instead of being converted to a Number.

$ws.Range("D2").Value = "27.457.58"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "1.832.04"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  -3.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.59"
$ws.Range("E5").Value = "  -2.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -2.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4305"
$ws.Range("E7").Value = "  -2.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3710"
$ws.Range("E8").Value = "  -2.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07267"
$ws.Range("E9").Value = "  -2.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8690"
$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.21"
$ws.Range("E11").Value = "  -2.06%  "

$ws.Range("D12").Value = "1.840.77"
$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.694"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.370"
$ws.Range("E14").Value = "  -3.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07087"
$ws.Range("E15").Value = "  -1.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.02"
$ws.Range("E16").Value = "  +2.25%  "

$ws.Range("E17").Value = "  -3.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008928"
$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.29"
$ws.Range("E20").Value = "  -1.94%  "

$ws.Range("D21").Value = "27.458.77"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.177"
$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.89"
$ws.Range("E23").Value = "  -3.53%  "

$ws.Range("D24").Value = "2.060.36"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.017"
$ws.Range("E25").Value = "  -2.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.55"
$ws.Range("E26").Value = "  -3.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  -1.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.154"
$ws.Range("E28").Value = "  +7.02%  "

$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.56"
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08851"
$ws.Range("E31").Value = "  -2.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.211"
$ws.Range("E32").Value = "  -0.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7707"
$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.511"
$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.901"
$ws.Range("E35").Value = "  -4.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  -3.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.124"
$ws.Range("E37").Value = "  -2.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01967"
$ws.Range("E38").Value = "  -0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05295"
$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.190"
$ws.Range("E40").Value = "  +3.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.891"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5091"
$ws.Range("E43").Value = "  -2.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.714"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.62"
$ws.Range("E46").Value = "  -4.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4745"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06426"
$ws.Range("E48").Value = "  -2.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("E49").Value = "  -3.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.676"
$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.835"
$ws.Range("E51").Value = "  -3.16%  "
